$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Quantum Computing: Unveiling a New Era of Technology" "A Dive into the Realm of Physics: Understanding the Fabric of Our Universe"

# Author name (merges "Dr" + "." + " Alex Hamilton" into one run)
Replace-Text "Dr. Alex Hamilton" "Grant Carter"

# Email user part
Replace-Text "alex" "grantcarter@phys"

# Email domain part (merges "hamilton@quantumresearch" + "." + "org" into one run)
Replace-Text "hamilton@quantumresearch.org" "edu"

# Paragraph 1, sentence 1 (keep the following "." run intact)
Replace-Text "Quantum mechanics, the enigmatic realm of the subatomic world, has long captivated the imagination of scientists and researchers" "Physics, a captivating scientific odyssey, unravels the secrets of our cosmos; it's the science that seeks to comprehend the underlying workings of the universe and its intricate composition"

# Paragraph 1, sentences 2-5 (merge into one run)
Replace-Text " As we delve deeper into the intricacies of this realm, we encounter phenomena that defy our conventional understanding of physics. Quantum computing, a rapidly emerging field, harnesses these enigmatic properties to construct computers capable of solving complex problems that are currently intractable for classical computers. By exploiting superposition and entanglement, quantum bits, or qubits, can exist in multiple states simultaneously and become inextricably linked to each other, unlocking unprecedented levels of computation. This transformative technology promises to revolutionize diverse industries, from cryptography to finance, medicine, and materials science" " It's a fascinating and intellectually challenging pursuit that investigates the fundamental principles that govern the behavior of the universe, from the tiniest subatomic particles to the grandest stellar formations"

# Paragraph 1, after first <br/><br/>
Replace-Text "The potential applications of quantum computing are as vast as they are intriguing" "The study of physics is akin to deciphering an enigmatic puzzle, where each discovery unveils another layer of the cosmic symphony"

Replace-Text " With the ability to solve previously unsolvable problems, quantum computers could accelerate drug discovery by simulating molecular interactions at lightning speed" " With each new revelation, we unravel the tapestry of the cosmos, revealing its hidden depths and illuminating the fundamental forces that mold its structure"

Replace-Text " They hold the potential to revolutionize cryptography by rendering current encryption methods obsolete, sparking a race to develop quantum-safe algorithms. In the realm of finance, quantum computing could optimize complex investment strategies and predict market trends with unprecedented accuracy. Furthermore, materials science could witness breakthroughs in the design of innovative materials with tailored properties, paving the way for advancements in energy storage, electronics, and medical devices" " Physics offers us a lens through which to understand the universe's complexities, from the grandeur of galaxies to the intimate world of atoms and subatomic particles"

# Paragraph 1, after second <br/><br/>
Replace-Text "As we stand at the cusp of a quantum revolution, there exists an urgent need to foster collaboration and knowledge sharing among researchers, industry leaders, and policymakers" "With curiosity as our compass and intellect as our guide, we embark on this extraordinary voyage into the realm of physics, seeking to unravel the mysteries of matter, energy, and the cosmos that encompasses us"

Replace-Text " International cooperation is paramount in addressing the challenges and harnessing the immense potential of this transformative technology. By establishing global initiatives, we can accelerate the development of quantum-ready workforce, fund cutting-edge research, and create an ecosystem that nurtures innovation. Moreover, it is essential to implement robust security measures to mitigate the risks associated with quantum computing, ensuring that this technology is used for the betterment of society and not for malicious purposes" " The quest for knowledge and the excitement of discovery fuel this journey as we explore the profound interconnectedness of all phenomena, from the macroscopic to the microscopic"

# Summary paragraph
Replace-Text "Quantum computing, an emerging field at the convergence of physics and computer science, has the potential to revolutionize various industries" "Physics, the study of the fundamental principles governing the universe, is an intellectually captivating pursuit that offers a deeper understanding of the cosmos and its components"

Replace-Text " By exploiting the principles of superposition and entanglement, quantum computers possess unprecedented computational capabilities that can solve complex problems intractable for classical computers" " We explore the fabric of reality through this scientific discipline, investigating the mysteries of matter, energy, and the cosmos, from the vastness of galaxies to the subatomic realm"

Replace-Text " The applications of quantum computing span a wide range of fields, including cryptography, finance, medicine, and materials science. However, to fully harness this transformative technology, international collaboration, investment in research, and the development of a quantum-ready workforce are crucial. As we embark on this quantum journey, we must also prioritize security measures to mitigate potential risks and ensure the responsible use of this technology. Embracing quantum computing has the power to usher in a new era of technological advancement and societal progress" " The journey into physics challenges us to decipher the intricate tapestry of the universe, unraveling its secrets and uncovering the underlying principles that shape our reality"

# Append trailing empty paragraph at the very end of the document body
$lastRange = $d.Content
$lastRange.Collapse(0)
$lastRange.InsertParagraphAfter()

Write-Host "done"
